# Zugdidi Municipality area sheet: drop the 1989/2002 census columns and the
# "(according to the population census data)" note, keeping just the 2014
# area figure, and restore the simpler (pre-coauthoring) row layout with
# spacer rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 1989 (B) and 2002 (C) columns -- only the 2014 figure (old
# column D) is kept, and it slides left into column B.
$ws.Range("B:C").EntireColumn.Delete()

# Remove the "(according to the population census data)" note row; the rows
# below shift up one (the blank spacer row becomes row 2, "(sq. km)" becomes
# row 3, the year header becomes row 4, and the Area figure becomes row 5).
$ws.Range("2:2").EntireRow.Delete()

# Restore the uniform 20.1pt row height used across the whole table
# (including a few trailing blank spacer rows) instead of the mixed
# 15 / 12.75 / 13.5 heights left over from the newer layout.
for ($i = 1; $i -le 8; $i++) {
    $ws.Rows.Item($i).RowHeight = 20.1
}
